$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 398. This pushes the
# existing row 398 (and everything below it, through row 484) down by
# one row, so the former row 484 becomes row 485.
$ws.Rows.Item(398).Insert()

# Populate the newly inserted row 398 with the new weekly price record,
# using the same fixed/template values (market, region, category, etc.)
# shared by all the other rows in this block, and the new date/price
# figures for this entry.
$ws.Range("A398").Value = 8
$ws.Range("B398").Value = "Terminal La Palmera de La Serena"
$ws.Range("C398").Value = "Coquimbo"
$ws.Range("D398").Value = 45211
$ws.Range("E398").Value = 4
$ws.Range("F398").Value = 100112012
$ws.Range("G398").Value = "Espinaca"
$ws.Range("H398").Value = "Sin especificar"
$ws.Range("I398").Value = "Primera"
$ws.Range("J398").Value = 1600
$ws.Range("K398").Value = 450
$ws.Range("L398").Value = 500
$ws.Range("M398").Value = 475
$ws.Range("N398").Value = "$/atado 300 a 500 gramos"
$ws.Range("O398").Value = "Provincia del Elquí"
$ws.Range("P398").Value = 950
$ws.Range("Q398").Value = 0.5
$ws.Range("R398").Value = "Hortaliza"
